$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-08-24 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-08-25 Monday", 2) | Out-Null
$d.Content.Find.Execute("37+6=43", $true, $false, $false, $false, $false, $true, 1, $false, "99+0=99", 2) | Out-Null
$d.Content.Find.Execute("69-48=21", $true, $false, $false, $false, $false, $true, 1, $false, "93-63=30", 2) | Out-Null
$d.Content.Find.Execute("14+40=54", $true, $false, $false, $false, $false, $true, 1, $false, "75-68=7", 2) | Out-Null
$d.Content.Find.Execute("23+68=91", $true, $false, $false, $false, $false, $true, 1, $false, "36+36=72", 2) | Out-Null
$d.Content.Find.Execute("26+26=52", $true, $false, $false, $false, $false, $true, 1, $false, "94-3=91", 2) | Out-Null
$d.Content.Find.Execute("49+26=75", $true, $false, $false, $false, $false, $true, 1, $false, "58-37=21", 2) | Out-Null
$d.Content.Find.Execute("35+59=94", $true, $false, $false, $false, $false, $true, 1, $false, "83+16=99", 2) | Out-Null
$d.Content.Find.Execute("5+76=81", $true, $false, $false, $false, $false, $true, 1, $false, "2+17=19", 2) | Out-Null
$d.Content.Find.Execute("67-25=42", $true, $false, $false, $false, $false, $true, 1, $false, "55-18=37", 2) | Out-Null
$d.Content.Find.Execute("97-7=90", $true, $false, $false, $false, $false, $true, 1, $false, "51-2=49", 2) | Out-Null
$d.Content.Find.Execute("90+4=94", $true, $false, $false, $false, $false, $true, 1, $false, "17+53=70", 2) | Out-Null
$d.Content.Find.Execute("33+8=41", $true, $false, $false, $false, $false, $true, 1, $false, "25+8=33", 2) | Out-Null
$d.Content.Find.Execute("52-50=2", $true, $false, $false, $false, $false, $true, 1, $false, "89-41=48", 2) | Out-Null
$d.Content.Find.Execute("89-82=7", $true, $false, $false, $false, $false, $true, 1, $false, "27+56=83", 2) | Out-Null
$d.Content.Find.Execute("48-36=12", $true, $false, $false, $false, $false, $true, 1, $false, "63-35=28", 2) | Out-Null
$d.Content.Find.Execute("37-17=20", $true, $false, $false, $false, $false, $true, 1, $false, "71+8=79", 2) | Out-Null
$d.Content.Find.Execute("44-33=11", $true, $false, $false, $false, $false, $true, 1, $false, "16+51=67", 2) | Out-Null
$d.Content.Find.Execute("92+4=96", $true, $false, $false, $false, $false, $true, 1, $false, "19+11=30", 2) | Out-Null
$d.Content.Find.Execute("87-46=41", $true, $false, $false, $false, $false, $true, 1, $false, "82-77=5", 2) | Out-Null
$d.Content.Find.Execute("84+14=98", $true, $false, $false, $false, $false, $true, 1, $false, "74-0=74", 2) | Out-Null
$d.Content.Find.Execute("3+20=23", $true, $false, $false, $false, $false, $true, 1, $false, "24+9=33", 2) | Out-Null
$d.Content.Find.Execute("14+36=50", $true, $false, $false, $false, $false, $true, 1, $false, "6+81=87", 2) | Out-Null
$d.Content.Find.Execute("13+56=69", $true, $false, $false, $false, $false, $true, 1, $false, "82-81=1", 2) | Out-Null
$d.Content.Find.Execute("52-51=1", $true, $false, $false, $false, $false, $true, 1, $false, "55-6=49", 2) | Out-Null
$d.Content.Find.Execute("11+25=36", $true, $false, $false, $false, $false, $true, 1, $false, "33+5=38", 2) | Out-Null
$d.Content.Find.Execute("6+17=23", $true, $false, $false, $false, $false, $true, 1, $false, "87-62=25", 2) | Out-Null
$d.Content.Find.Execute("14+12=26", $true, $false, $false, $false, $false, $true, 1, $false, "38-38=0", 2) | Out-Null
$d.Content.Find.Execute("40-9=31", $true, $false, $false, $false, $false, $true, 1, $false, "36+50=86", 2) | Out-Null
$d.Content.Find.Execute("24+36=60", $true, $false, $false, $false, $false, $true, 1, $false, "66-41=25", 2) | Out-Null
$d.Content.Find.Execute("99-95=4", $true, $false, $false, $false, $false, $true, 1, $false, "15-4=11", 2) | Out-Null
$d.Content.Find.Execute("13+42=55", $true, $false, $false, $false, $false, $true, 1, $false, "95-6=89", 2) | Out-Null
$d.Content.Find.Execute("52+33=85", $true, $false, $false, $false, $false, $true, 1, $false, "10+76=86", 2) | Out-Null
$d.Content.Find.Execute("99-22=77", $true, $false, $false, $false, $false, $true, 1, $false, "43-31=12", 2) | Out-Null
$d.Content.Find.Execute("11+40=51", $true, $false, $false, $false, $false, $true, 1, $false, "67-9=58", 2) | Out-Null
$d.Content.Find.Execute("7+0=7", $true, $false, $false, $false, $false, $true, 1, $false, "13-3=10", 2) | Out-Null
$d.Content.Find.Execute("7+61=68", $true, $false, $false, $false, $false, $true, 1, $false, "97-27=70", 2) | Out-Null
$d.Content.Find.Execute("94-29=65", $true, $false, $false, $false, $false, $true, 1, $false, "62-38=24", 2) | Out-Null
$d.Content.Find.Execute("81-16=65", $true, $false, $false, $false, $false, $true, 1, $false, "25-10=15", 2) | Out-Null
$d.Content.Find.Execute("66-30=36", $true, $false, $false, $false, $false, $true, 1, $false, "71-67=4", 2) | Out-Null
$d.Content.Find.Execute("70-25=45", $true, $false, $false, $false, $false, $true, 1, $false, "34+34=68", 2) | Out-Null
$d.Content.Find.Execute("33-22=11", $true, $false, $false, $false, $false, $true, 1, $false, "18+12=30", 2) | Out-Null
$d.Content.Find.Execute("96-40=56", $true, $false, $false, $false, $false, $true, 1, $false, "50+15=65", 2) | Out-Null
$d.Content.Find.Execute("22+19=41", $true, $false, $false, $false, $false, $true, 1, $false, "44-30=14", 2) | Out-Null
$d.Content.Find.Execute("99-96=3", $true, $false, $false, $false, $false, $true, 1, $false, "71-62=9", 2) | Out-Null
$d.Content.Find.Execute("79-36=43", $true, $false, $false, $false, $false, $true, 1, $false, "95-86=9", 2) | Out-Null
$d.Content.Find.Execute("46-26=20", $true, $false, $false, $false, $false, $true, 1, $false, "36-21=15", 2) | Out-Null
$d.Content.Find.Execute("20+70=90", $true, $false, $false, $false, $false, $true, 1, $false, "54-52=2", 2) | Out-Null
$d.Content.Find.Execute("29+70=99", $true, $false, $false, $false, $false, $true, 1, $false, "4+91=95", 2) | Out-Null
$d.Content.Find.Execute("90-11=79", $true, $false, $false, $false, $false, $true, 1, $false, "8-7=1", 2) | Out-Null
$d.Content.Find.Execute("76+5=81", $true, $false, $false, $false, $false, $true, 1, $false, "12+35=47", 2) | Out-Null
$d.Content.Find.Execute("61+0=61", $true, $false, $false, $false, $false, $true, 1, $false, "28-23=5", 2) | Out-Null
$d.Content.Find.Execute("32+56=88", $true, $false, $false, $false, $false, $true, 1, $false, "63+23=86", 2) | Out-Null
$d.Content.Find.Execute("26+62=88", $true, $false, $false, $false, $false, $true, 1, $false, "31-16=15", 2) | Out-Null
$d.Content.Find.Execute("17+77=94", $true, $false, $false, $false, $false, $true, 1, $false, "16+7=23", 2) | Out-Null
$d.Content.Find.Execute("45-36=9", $true, $false, $false, $false, $false, $true, 1, $false, "58-25=33", 2) | Out-Null
$d.Content.Find.Execute("73-72=1", $true, $false, $false, $false, $false, $true, 1, $false, "61-11=50", 2) | Out-Null
$d.Content.Find.Execute("96-18=78", $true, $false, $false, $false, $false, $true, 1, $false, "21+64=85", 2) | Out-Null
$d.Content.Find.Execute("44+19=63", $true, $false, $false, $false, $false, $true, 1, $false, "79-78=1", 2) | Out-Null
$d.Content.Find.Execute("7+10=17", $true, $false, $false, $false, $false, $true, 1, $false, "47-13=34", 2) | Out-Null
$d.Content.Find.Execute("15+10=25", $true, $false, $false, $false, $false, $true, 1, $false, "27+41=68", 2) | Out-Null
$d.Content.Find.Execute("47+2=49", $true, $false, $false, $false, $false, $true, 1, $false, "99-83=16", 2) | Out-Null
$d.Content.Find.Execute("48-44=4", $true, $false, $false, $false, $false, $true, 1, $false, "41-23=18", 2) | Out-Null
$d.Content.Find.Execute("28-13=15", $true, $false, $false, $false, $false, $true, 1, $false, "90-85=5", 2) | Out-Null
$d.Content.Find.Execute("67-21=46", $true, $false, $false, $false, $false, $true, 1, $false, "88-18=70", 2) | Out-Null
$d.Content.Find.Execute("2+85=87", $true, $false, $false, $false, $false, $true, 1, $false, "19-4=15", 2) | Out-Null
$d.Content.Find.Execute("79-29=50", $true, $false, $false, $false, $false, $true, 1, $false, "31+57=88", 2) | Out-Null
$d.Content.Find.Execute("82-31=51", $true, $false, $false, $false, $false, $true, 1, $false, "78-38=40", 2) | Out-Null
$d.Content.Find.Execute("64-19=45", $true, $false, $false, $false, $false, $true, 1, $false, "24+45=69", 2) | Out-Null
$d.Content.Find.Execute("1+59=60", $true, $false, $false, $false, $false, $true, 1, $false, "74+21=95", 2) | Out-Null
$d.Content.Find.Execute("17+57=74", $true, $false, $false, $false, $false, $true, 1, $false, "95-62=33", 2) | Out-Null
$d.Content.Find.Execute("56+12=68", $true, $false, $false, $false, $false, $true, 1, $false, "25-9=16", 2) | Out-Null
$d.Content.Find.Execute("39+54=93", $true, $false, $false, $false, $false, $true, 1, $false, "79-70=9", 2) | Out-Null
$d.Content.Find.Execute("39+44=83", $true, $false, $false, $false, $false, $true, 1, $false, "26+9=35", 2) | Out-Null
$d.Content.Find.Execute("54+3=57", $true, $false, $false, $false, $false, $true, 1, $false, "45-4=41", 2) | Out-Null
$d.Content.Find.Execute("15+81=96", $true, $false, $false, $false, $false, $true, 1, $false, "70-54=16", 2) | Out-Null
$d.Content.Find.Execute("56+28=84", $true, $false, $false, $false, $false, $true, 1, $false, "40+48=88", 2) | Out-Null
$d.Content.Find.Execute("94-45=49", $true, $false, $false, $false, $false, $true, 1, $false, "91-79=12", 2) | Out-Null
$d.Content.Find.Execute("19+66=85", $true, $false, $false, $false, $false, $true, 1, $false, "28+37=65", 2) | Out-Null
$d.Content.Find.Execute("15+38=53", $true, $false, $false, $false, $false, $true, 1, $false, "87-25=62", 2) | Out-Null
$d.Content.Find.Execute("50+26=76", $true, $false, $false, $false, $false, $true, 1, $false, "53-47=6", 2) | Out-Null
$d.Content.Find.Execute("30+18=48", $true, $false, $false, $false, $false, $true, 1, $false, "87-57=30", 2) | Out-Null
$d.Content.Find.Execute("57-17=40", $true, $false, $false, $false, $false, $true, 1, $false, "45+49=94", 2) | Out-Null
$d.Content.Find.Execute("19+45=64", $true, $false, $false, $false, $false, $true, 1, $false, "23+53=76", 2) | Out-Null
$d.Content.Find.Execute("93-10=83", $true, $false, $false, $false, $false, $true, 1, $false, "65-34=31", 2) | Out-Null
$d.Content.Find.Execute("97-10=87", $true, $false, $false, $false, $false, $true, 1, $false, "48-28=20", 2) | Out-Null
$d.Content.Find.Execute("12-9=3", $true, $false, $false, $false, $false, $true, 1, $false, "17+36=53", 2) | Out-Null
$d.Content.Find.Execute("29+60=89", $true, $false, $false, $false, $false, $true, 1, $false, "58+6=64", 2) | Out-Null
$d.Content.Find.Execute("45+40=85", $true, $false, $false, $false, $false, $true, 1, $false, "32+15=47", 2) | Out-Null
$d.Content.Find.Execute("0+47=47", $true, $false, $false, $false, $false, $true, 1, $false, "57+7=64", 2) | Out-Null
$d.Content.Find.Execute("9+41=50", $true, $false, $false, $false, $false, $true, 1, $false, "49+14=63", 2) | Out-Null
$d.Content.Find.Execute("12+21=33", $true, $false, $false, $false, $false, $true, 1, $false, "2+25=27", 2) | Out-Null
$d.Content.Find.Execute("49-26=23", $true, $false, $false, $false, $false, $true, 1, $false, "47+39=86", 2) | Out-Null
$d.Content.Find.Execute("45+53=98", $true, $false, $false, $false, $false, $true, 1, $false, "31-12=19", 2) | Out-Null
$d.Content.Find.Execute("17-0=17", $true, $false, $false, $false, $false, $true, 1, $false, "76+20=96", 2) | Out-Null
$d.Content.Find.Execute("58+22=80", $true, $false, $false, $false, $false, $true, 1, $false, "87+12=99", 2) | Out-Null
$d.Content.Find.Execute("63+19=82", $true, $false, $false, $false, $false, $true, 1, $false, "28+33=61", 2) | Out-Null
$d.Content.Find.Execute("81-74=7", $true, $false, $false, $false, $false, $true, 1, $false, "8+35=43", 2) | Out-Null
$d.Content.Find.Execute("14+29=43", $true, $false, $false, $false, $false, $true, 1, $false, "12+58=70", 2) | Out-Null
$d.Content.Find.Execute("13+83=96", $true, $false, $false, $false, $false, $true, 1, $false, "31-11=20", 2) | Out-Null
$d.Content.Find.Execute("39+22=61", $true, $false, $false, $false, $false, $true, 1, $false, "4+13=17", 2) | Out-Null
